$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 687.6667
$ws.Range("I8").Value = 63
$ws.Range("K8").Value = 189
$ws.Range("M8").Value = -50

$ws.Range("H17").Value = 5533228
$ws.Range("J17").Value = 5533228
$ws.Range("L17").Value = 16599684
$ws.Range("N17").Value = -16600020

$ws.Range("H19").Value = 521.6429000000001
$ws.Range("I19").Value = 428.625
$ws.Range("J19").Value = 645.6667
$ws.Range("K19").Value = 428.625
$ws.Range("L19").Value = 645.6667
$ws.Range("M19").Value = -253.625
$ws.Range("N19").Value = -995.6667

$ws.Range("H21").Value = 65009.5
$ws.Range("I21").Value = 80019
$ws.Range("J21").Value = 50000
$ws.Range("K21").Value = 80019
$ws.Range("L21").Value = 50000
$ws.Range("M21").Value = -79551
$ws.Range("N21").Value = -50936

$ws.Range("H23").Value = 65009.5
$ws.Range("I23").Value = 80019
$ws.Range("J23").Value = 50000
$ws.Range("K23").Value = 80019
$ws.Range("L23").Value = 50000
$ws.Range("M23").Value = -79785
$ws.Range("N23").Value = -50468

$ws.Range("H64").Value = 4605.909
$ws.Range("I64").Value = 5239.2856
$ws.Range("J64").Value = 3497.5
$ws.Range("K64").Value = 5239.2856
$ws.Range("L64").Value = 3497.5
$ws.Range("M64").Value = -4991.2856
$ws.Range("N64").Value = -3993.5

$ws.Range("H67").Value = 4605.909
$ws.Range("I67").Value = 5239.2856
$ws.Range("J67").Value = 3497.5
$ws.Range("K67").Value = 5239.2856
$ws.Range("L67").Value = 3497.5
$ws.Range("M67").Value = -4381.2856
$ws.Range("N67").Value = -5213.5

$ws.Range("H74").Value = 3780.6
$ws.Range("I74").Value = 3742.4167
$ws.Range("J74").Value = 3933.3333
$ws.Range("K74").Value = 3742.4167
$ws.Range("L74").Value = 3933.3333
$ws.Range("M74").Value = -2806.4167
$ws.Range("N74").Value = -5805.3333

$ws.Range("H77").Value = 3780.6
$ws.Range("I77").Value = 3742.4167
$ws.Range("J77").Value = 3933.3333
$ws.Range("K77").Value = 18712.0835
$ws.Range("L77").Value = 19666.6665
$ws.Range("M77").Value = -14032.0835
$ws.Range("N77").Value = -29026.6665

$ws.Range("H106").Value = 10756706
$ws.Range("I106").Value = 11498392
$ws.Range("J106").Value = 2253
$ws.Range("K106").Value = 11498392
$ws.Range("L106").Value = 2253
$ws.Range("M106").Value = -11497761
$ws.Range("N106").Value = -3515

$ws.Range("H132").Value = 2059.9285
$ws.Range("I132").Value = 1489.8108
$ws.Range("J132").Value = 6278.8
$ws.Range("K132").Value = 4469.4324
$ws.Range("L132").Value = 18836.4
$ws.Range("M132").Value = -1939.4324
$ws.Range("N132").Value = -23896.4

$ws.Range("H133").Value = 51187.5
$ws.Range("J133").Value = 51187.5
$ws.Range("L133").Value = 51187.5
$ws.Range("N133").Value = -61307.5

$ws.Range("H135").Value = 2507.04
$ws.Range("I135").Value = 2299
$ws.Range("J135").Value = 3165.8333
$ws.Range("K135").Value = 20691
$ws.Range("L135").Value = 28492.4997
$ws.Range("M135").Value = -18156
$ws.Range("N135").Value = -33562.4997

$ws.Range("H138").Value = 2325.1953
$ws.Range("I138").Value = 1471.2916
$ws.Range("J138").Value = 2650.492
$ws.Range("K138").Value = 4413.8748
$ws.Range("L138").Value = 7951.476000000001
$ws.Range("M138").Value = 726.1252000000004
$ws.Range("N138").Value = -18231.476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 5946.154
$ws.Range("I6").Value = 1666.6666
$ws.Range("J6").Value = 7230
$ws.Range("K6").Value = 1666.6666
$ws.Range("L6").Value = 7230
$ws.Range("M6").Value = -1493.6666
$ws.Range("N6").Value = -7576

$ws.Range("H32").Value = 3185.763
$ws.Range("I32").Value = 3185.763
$ws.Range("K32").Value = 3185.763
$ws.Range("M32").Value = -2898.763

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 52500
$ws.Range("J115").Value = 60000
$ws.Range("L115").Value = 60000
$ws.Range("N115").Value = -63134

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 160.90909
$ws.Range("I22").Value = 140
$ws.Range("J22").Value = 197.5
$ws.Range("K22").Value = 140
$ws.Range("L22").Value = 197.5
$ws.Range("M22").Value = 210
$ws.Range("N22").Value = -897.5

$ws.Range("H31").Value = 4577.706
$ws.Range("I31").Value = 1936.7778
$ws.Range("J31").Value = 9744.739
$ws.Range("K31").Value = 1936.7778
$ws.Range("L31").Value = 9744.739
$ws.Range("M31").Value = -1641.7778
$ws.Range("N31").Value = -10334.739

$ws.Range("H34").Value = 4577.706
$ws.Range("I34").Value = 1936.7778
$ws.Range("J34").Value = 9744.739
$ws.Range("K34").Value = 1936.7778
$ws.Range("L34").Value = 9744.739
$ws.Range("M34").Value = -1734.7778
$ws.Range("N34").Value = -10148.739

$ws.Range("H58").Value = 1923.8889
$ws.Range("I58").Value = 1541.9445
$ws.Range("J58").Value = 2178.5186
$ws.Range("K58").Value = 1541.9445
$ws.Range("L58").Value = 2178.5186
$ws.Range("M58").Value = -1338.9445
$ws.Range("N58").Value = -2584.5186

$ws.Range("H99").Value = 17080.6
$ws.Range("I99").Value = 27135.334
$ws.Range("J99").Value = 1998.5
$ws.Range("K99").Value = 27135.334
$ws.Range("L99").Value = 1998.5
$ws.Range("M99").Value = -25637.334
$ws.Range("N99").Value = -4994.5

$ws.Range("H126").Value = 17080.6
$ws.Range("I126").Value = 27135.334
$ws.Range("J126").Value = 1998.5
$ws.Range("K126").Value = 81406.00199999999
$ws.Range("L126").Value = 5995.5
$ws.Range("M126").Value = -78936.00199999999
$ws.Range("N126").Value = -10935.5

$ws.Range("H132").Value = 2210.087
$ws.Range("I132").Value = 1137.8462
$ws.Range("J132").Value = 3604
$ws.Range("K132").Value = 3413.5386
$ws.Range("L132").Value = 10812
$ws.Range("M132").Value = -883.5385999999999
$ws.Range("N132").Value = -15872

$ws.Range("H136").Value = 1923.8889
$ws.Range("I136").Value = 1541.9445
$ws.Range("J136").Value = 2178.5186
$ws.Range("K136").Value = 4625.833500000001
$ws.Range("L136").Value = 6535.5558
$ws.Range("M136").Value = -2075.833500000001
$ws.Range("N136").Value = -11635.5558

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2022.4
$ws.Range("I2").Value = 2518
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 15108
$ws.Range("L2").Value = 240
$ws.Range("M2").Value = -14995
$ws.Range("N2").Value = -466

$ws.Range("H4").Value = 598730.9
$ws.Range("I4").Value = 2386923.8
$ws.Range("J4").Value = 2666.611
$ws.Range("K4").Value = 7160771.399999999
$ws.Range("L4").Value = 7999.833
$ws.Range("M4").Value = -7160659.399999999
$ws.Range("N4").Value = -8223.832999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10636454
$ws.Range("I11").Value = 10143000
$ws.Range("J11").Value = 11500000
$ws.Range("K11").Value = 10143000
$ws.Range("L11").Value = 11500000
$ws.Range("M11").Value = -10142861
$ws.Range("N11").Value = -11500278

$ws.Range("H102").Value = 1871.6086
$ws.Range("I102").Value = 1792.35
$ws.Range("J102").Value = 2400
$ws.Range("K102").Value = 1792.35
$ws.Range("L102").Value = 2400
$ws.Range("M102").Value = -170.3499999999999
$ws.Range("N102").Value = -5644

$ws.Range("H126").Value = 7382.1387
$ws.Range("I126").Value = 8835.444
$ws.Range("J126").Value = 3022.2222
$ws.Range("K126").Value = 26506.332
$ws.Range("L126").Value = 9066.6666
$ws.Range("M126").Value = -24036.332
$ws.Range("N126").Value = -14006.6666

$ws.Range("H132").Value = 5803.5527
$ws.Range("I132").Value = 6474.407
$ws.Range("J132").Value = 4156.909
$ws.Range("K132").Value = 19423.221
$ws.Range("L132").Value = 12470.727
$ws.Range("M132").Value = -16893.221
$ws.Range("N132").Value = -17530.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 90912140
$ws.Range("I40").Value = 100003010
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 100003010
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -100002874
$ws.Range("N40").Value = -3772

$ws.Range("H61").Value = 1311.4
$ws.Range("I61").Value = 1233.0526
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 1233.0526
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -1031.0526
$ws.Range("N61").Value = -3204

$ws.Range("H68").Value = 2400.7
$ws.Range("I68").Value = 2138
$ws.Range("K68").Value = 2138
$ws.Range("M68").Value = -1389

$ws.Range("H71").Value = 2400.7
$ws.Range("I71").Value = 2138
$ws.Range("K71").Value = 10690
$ws.Range("M71").Value = -6946

$ws.Range("H82").Value = 229400
$ws.Range("I82").Value = 9000
$ws.Range("J82").Value = 376333.34
$ws.Range("K82").Value = 9000
$ws.Range("L82").Value = 376333.34
$ws.Range("M82").Value = -8639
$ws.Range("N82").Value = -377055.34

$ws.Range("H85").Value = 229400
$ws.Range("I85").Value = 9000
$ws.Range("J85").Value = 376333.34
$ws.Range("K85").Value = 9000
$ws.Range("L85").Value = 376333.34
$ws.Range("M85").Value = -7752
$ws.Range("N85").Value = -378829.34

$ws.Range("H100").Value = 1335.3846
$ws.Range("I100").Value = 1335.3846
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1335.3846
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -794.3846000000001
$ws.Range("N100").ClearContents()

$ws.Range("H113").Value = 1311.4
$ws.Range("I113").Value = 1233.0526
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 1233.0526
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = 936.9474
$ws.Range("N113").Value = -7140

$ws.Range("H136").Value = 6118.925
$ws.Range("I136").Value = 6757.0435
$ws.Range("J136").Value = 5255.5884
$ws.Range("K136").Value = 20271.1305
$ws.Range("L136").Value = 15766.7652
$ws.Range("M136").Value = -17721.1305
$ws.Range("N136").Value = -20866.7652

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 7220
$ws.Range("J22").Value = 7220
$ws.Range("L22").Value = 7220
$ws.Range("N22").Value = -7806

$ws.Range("H107").Value = 40000668
$ws.Range("I107").Value = 90909600
$ws.Range("J107").Value = 791.1429000000001
$ws.Range("K107").Value = 272728800
$ws.Range("L107").Value = 2373.4287
$ws.Range("M107").Value = -272726880
$ws.Range("N107").Value = -6213.4287

$ws.Range("H136").Value = 2204.05
$ws.Range("I136").Value = 2163.4167
$ws.Range("J136").Value = 2265
$ws.Range("K136").Value = 6490.250100000001
$ws.Range("L136").Value = 6795
$ws.Range("M136").Value = -3940.250100000001
$ws.Range("N136").Value = -11895
